{"js": "// Replacements derived from the diff: each [before, after] pair is a\n// unique, exact text value that appears exactly once in the document\n// body (the date paragraph + 25 distinct table-cell \"A\u00f7B=C, D\" strings).\nconst replacements = [\n  [\"2024-02-25 Sunday\", \"2024-02-26 Monday\"],\n  [\"201\u00f72=100, 1\", \"399\u00f75=79, 4\"],\n  [\"929\u00f75=185, 4\", \"881\u00f72=440, 1\"],\n  [\"447\u00f73=149, 0\", \"451\u00f77=64, 3\"],\n  [\"922\u00f76=153, 4\", \"310\u00f79=34, 4\"],\n  [\"448\u00f79=49, 7\", \"354\u00f74=88, 2\"],\n  [\"508\u00f78=63, 4\", \"255\u00f72=127, 1\"],\n  [\"703\u00f75=140, 3\", \"931\u00f76=155, 1\"],\n  [\"873\u00f77=124, 5\", \"324\u00f77=46, 2\"],\n  [\"170\u00f79=18, 8\", \"791\u00f77=113, 0\"],\n  [\"818\u00f77=116, 6\", \"189\u00f76=31, 3\"],\n  [\"847\u00f78=105, 7\", \"564\u00f76=94, 0\"],\n  [\"956\u00f78=119, 4\", \"117\u00f78=14, 5\"],\n  [\"630\u00f76=105, 0\", \"524\u00f76=87, 2\"],\n  [\"880\u00f77=125, 5\", \"539\u00f75=107, 4\"],\n  [\"380\u00f72=190, 0\", \"902\u00f74=225, 2\"],\n  [\"505\u00f77=72, 1\", \"234\u00f76=39, 0\"],\n  [\"835\u00f78=104, 3\", \"562\u00f74=140, 2\"],\n  [\"270\u00f73=90, 0\", \"215\u00f74=53, 3\"],\n  [\"287\u00f77=41, 0\", \"323\u00f76=53, 5\"],\n  [\"270\u00f75=54, 0\", \"869\u00f75=173, 4\"],\n  [\"473\u00f75=94, 3\", \"288\u00f73=96, 0\"],\n  [\"197\u00f76=32, 5\", \"217\u00f74=54, 1\"],\n  [\"201\u00f77=28, 5\", \"292\u00f79=32, 4\"],\n  [\"278\u00f74=69, 2\", \"538\u00f79=59, 7\"],\n  [\"905\u00f73=301, 2\", \"466\u00f76=77, 4\"],\n];\n\nfor (const [before, after] of replacements) {\n  const found = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items/text\");\n  await context.sync();\n\n  // Exact, unique text -> exactly one hit expected; replace it (and any\n  // further occurrences, defensively) in place so run formatting survives.\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-25 Sunday\", \"2024-02-26 Monday\"),\n    @(\"201\u00f72=100, 1\", \"399\u00f75=79, 4\"),\n    @(\"929\u00f75=185, 4\", \"881\u00f72=440, 1\"),\n    @(\"447\u00f73=149, 0\", \"451\u00f77=64, 3\"),\n    @(\"922\u00f76=153, 4\", \"310\u00f79=34, 4\"),\n    @(\"448\u00f79=49, 7\", \"354\u00f74=88, 2\"),\n    @(\"508\u00f78=63, 4\", \"255\u00f72=127, 1\"),\n    @(\"703\u00f75=140, 3\", \"931\u00f76=155, 1\"),\n    @(\"873\u00f77=124, 5\", \"324\u00f77=46, 2\"),\n    @(\"170\u00f79=18, 8\", \"791\u00f77=113, 0\"),\n    @(\"818\u00f77=116, 6\", \"189\u00f76=31, 3\"),\n    @(\"847\u00f78=105, 7\", \"564\u00f76=94, 0\"),\n    @(\"956\u00f78=119, 4\", \"117\u00f78=14, 5\"),\n    @(\"630\u00f76=105, 0\", \"524\u00f76=87, 2\"),\n    @(\"880\u00f77=125, 5\", \"539\u00f75=107, 4\"),\n    @(\"380\u00f72=190, 0\", \"902\u00f74=225, 2\"),\n    @(\"505\u00f77=72, 1\", \"234\u00f76=39, 0\"),\n    @(\"835\u00f78=104, 3\", \"562\u00f74=140, 2\"),\n    @(\"270\u00f73=90, 0\", \"215\u00f74=53, 3\"),\n    @(\"287\u00f77=41, 0\", \"323\u00f76=53, 5\"),\n    @(\"270\u00f75=54, 0\", \"869\u00f75=173, 4\"),\n    @(\"473\u00f75=94, 3\", \"288\u00f73=96, 0\"),\n    @(\"197\u00f76=32, 5\", \"217\u00f74=54, 1\"),\n    @(\"201\u00f77=28, 5\", \"292\u00f79=32, 4\"),\n    @(\"278\u00f74=69, 2\", \"538\u00f79=59, 7\"),\n    @(\"905\u00f73=301, 2\", \"466\u00f76=77, 4\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null  # wdReplaceAll = 2\n}\n"}
